$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Investissements")
$ws1.Range("B3").Value = "AnnÃ©e"
$ws1.Range("C3").Value = "Investissements spÃ©cifiques"
$ws1.Range("D3").Value = "Investissements intÃ©grÃ©s"
$ws1.Range("F3").Value = "Ã‰tudes"
$ws1.Range("G3").Value = "Total investissements et Ã©tudes"
$ws2 = $wb.Worksheets.Item("Dépenses 1")
$ws2.Range("B3").Value = "DÃ©penses courantes"
$ws2.Range("B5").Value = "Fonctionnement des Ã©quipements dÃ©diÃ©s Ã  la protection de l'environnement"
$ws2.Range("B6").Value = "Autres dÃ©penses courantes"
$ws2.Range("B7").Value = "Total dÃ©penses courantes"
$ws3 = $wb.Worksheets.Item("Dépenses 2")
$ws3.Range("B3").Value = "DÃ©penses courantes"
$ws3.Range("B5").Value = "   dont liÃ©es Ã  l'eau"
$ws3.Range("B6").Value = "   dont liÃ©es aux dÃ©chets"
$ws3.Range("B7").Value = "Fonctionnement des Ã©quipements dÃ©diÃ©s Ã  la protection de l'environnement"
$ws3.Range("B8").Value = "    coÃ»ts internes "
$ws3.Range("B11").Value = "Autres dÃ©penses courantes"
$ws3.Range("B12").Value = "Total dÃ©penses courantes"
